{"js": "// Replace \"elaboration\" with \"Construction\" in the revision-history\n// table's Description cell (\"elaboration iteration 4 draft 1\" ->\n// \"Construction iteration 4 draft 1\").\nconst body = context.document.body;\nconst results = body.search(\"elaboration iteration 4 draft 1\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(\"Construction iteration 4 draft 1\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Replace \"elaboration\" with \"Construction\" in the revision-history\n# table's Description cell (\"elaboration iteration 4 draft 1\" ->\n# \"Construction iteration 4 draft 1\").\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"elaboration iteration 4 draft 1\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"Construction iteration 4 draft 1\"\n$find.Execute(\n    $find.Text,      # FindText\n    $false,          # MatchCase\n    $false,          # MatchWholeWord\n    $false,          # MatchWildcards\n    $false,          # MatchSoundsLike\n    $false,          # MatchAllWordForms\n    $true,           # Forward\n    1,               # Wrap (wdFindContinue)\n    $false,          # Format\n    $find.Replacement.Text,  # ReplaceWith\n    2                # Replace (wdReplaceAll)\n)\n"}
